# Updates the date heading and every arithmetic-answer cell in the table to
# the values from the "answers-of-within100" regeneration (commit 596fc94).
# Each old value is unique in the document, so a simple MatchCase, non-wildcard
# Find/Replace (ReplaceAll) for each pair is safe and order-independent.

$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-06-10 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-06-11 Sunday", 2)
$d.Content.Find.Execute("3+5=8", $true, $false, $false, $false, $false, $true, 1, $false, "42-28=14", 2)
$d.Content.Find.Execute("80-39=41", $true, $false, $false, $false, $false, $true, 1, $false, "42-33=9", 2)
$d.Content.Find.Execute("58-12=46", $true, $false, $false, $false, $false, $true, 1, $false, "90-52=38", 2)
$d.Content.Find.Execute("10+39=49", $true, $false, $false, $false, $false, $true, 1, $false, "98-24=74", 2)
$d.Content.Find.Execute("14+41=55", $true, $false, $false, $false, $false, $true, 1, $false, "1+39=40", 2)
$d.Content.Find.Execute("95-36=59", $true, $false, $false, $false, $false, $true, 1, $false, "94-57=37", 2)
$d.Content.Find.Execute("47-29=18", $true, $false, $false, $false, $false, $true, 1, $false, "76-22=54", 2)
$d.Content.Find.Execute("96-16=80", $true, $false, $false, $false, $false, $true, 1, $false, "52-3=49", 2)
$d.Content.Find.Execute("55-20=35", $true, $false, $false, $false, $false, $true, 1, $false, "37-6=31", 2)
$d.Content.Find.Execute("30+41=71", $true, $false, $false, $false, $false, $true, 1, $false, "8+53=61", 2)
$d.Content.Find.Execute("18+56=74", $true, $false, $false, $false, $false, $true, 1, $false, "94-65=29", 2)
$d.Content.Find.Execute("47-13=34", $true, $false, $false, $false, $false, $true, 1, $false, "75-39=36", 2)
$d.Content.Find.Execute("26+26=52", $true, $false, $false, $false, $false, $true, 1, $false, "94-21=73", 2)
$d.Content.Find.Execute("95-43=52", $true, $false, $false, $false, $false, $true, 1, $false, "44+2=46", 2)
$d.Content.Find.Execute("9+12=21", $true, $false, $false, $false, $false, $true, 1, $false, "28+57=85", 2)
$d.Content.Find.Execute("16-7=9", $true, $false, $false, $false, $false, $true, 1, $false, "99-77=22", 2)
$d.Content.Find.Execute("91-29=62", $true, $false, $false, $false, $false, $true, 1, $false, "70+2=72", 2)
$d.Content.Find.Execute("34-19=15", $true, $false, $false, $false, $false, $true, 1, $false, "80-22=58", 2)
$d.Content.Find.Execute("43-31=12", $true, $false, $false, $false, $false, $true, 1, $false, "95-0=95", 2)
$d.Content.Find.Execute("69-53=16", $true, $false, $false, $false, $false, $true, 1, $false, "41-7=34", 2)
$d.Content.Find.Execute("50+37=87", $true, $false, $false, $false, $false, $true, 1, $false, "20+35=55", 2)
$d.Content.Find.Execute("6-3=3", $true, $false, $false, $false, $false, $true, 1, $false, "70-42=28", 2)
$d.Content.Find.Execute("69+15=84", $true, $false, $false, $false, $false, $true, 1, $false, "74-48=26", 2)
$d.Content.Find.Execute("62-53=9", $true, $false, $false, $false, $false, $true, 1, $false, "72-3=69", 2)
$d.Content.Find.Execute("56-51=5", $true, $false, $false, $false, $false, $true, 1, $false, "64-7=57", 2)
$d.Content.Find.Execute("45+1=46", $true, $false, $false, $false, $false, $true, 1, $false, "21+22=43", 2)
$d.Content.Find.Execute("96-49=47", $true, $false, $false, $false, $false, $true, 1, $false, "7+29=36", 2)
$d.Content.Find.Execute("49-48=1", $true, $false, $false, $false, $false, $true, 1, $false, "32-5=27", 2)
$d.Content.Find.Execute("82-23=59", $true, $false, $false, $false, $false, $true, 1, $false, "38+3=41", 2)
$d.Content.Find.Execute("76-48=28", $true, $false, $false, $false, $false, $true, 1, $false, "97-68=29", 2)
$d.Content.Find.Execute("29+26=55", $true, $false, $false, $false, $false, $true, 1, $false, "95-19=76", 2)
$d.Content.Find.Execute("91-38=53", $true, $false, $false, $false, $false, $true, 1, $false, "53+14=67", 2)
$d.Content.Find.Execute("90-7=83", $true, $false, $false, $false, $false, $true, 1, $false, "67-31=36", 2)
$d.Content.Find.Execute("61+26=87", $true, $false, $false, $false, $false, $true, 1, $false, "64-0=64", 2)
$d.Content.Find.Execute("34+0=34", $true, $false, $false, $false, $false, $true, 1, $false, "0+56=56", 2)
$d.Content.Find.Execute("12+85=97", $true, $false, $false, $false, $false, $true, 1, $false, "1+30=31", 2)
$d.Content.Find.Execute("26+48=74", $true, $false, $false, $false, $false, $true, 1, $false, "44+51=95", 2)
$d.Content.Find.Execute("22-7=15", $true, $false, $false, $false, $false, $true, 1, $false, "7+14=21", 2)
$d.Content.Find.Execute("3+61=64", $true, $false, $false, $false, $false, $true, 1, $false, "93+2=95", 2)
$d.Content.Find.Execute("88-69=19", $true, $false, $false, $false, $false, $true, 1, $false, "86-75=11", 2)
$d.Content.Find.Execute("19+41=60", $true, $false, $false, $false, $false, $true, 1, $false, "7+24=31", 2)
$d.Content.Find.Execute("42+43=85", $true, $false, $false, $false, $false, $true, 1, $false, "59-19=40", 2)
$d.Content.Find.Execute("91-22=69", $true, $false, $false, $false, $false, $true, 1, $false, "51-20=31", 2)
$d.Content.Find.Execute("72-24=48", $true, $false, $false, $false, $false, $true, 1, $false, "96-68=28", 2)
$d.Content.Find.Execute("13-4=9", $true, $false, $false, $false, $false, $true, 1, $false, "42-30=12", 2)
$d.Content.Find.Execute("80+17=97", $true, $false, $false, $false, $false, $true, 1, $false, "10+35=45", 2)
$d.Content.Find.Execute("49+24=73", $true, $false, $false, $false, $false, $true, 1, $false, "28-19=9", 2)
$d.Content.Find.Execute("95-31=64", $true, $false, $false, $false, $false, $true, 1, $false, "36-24=12", 2)
$d.Content.Find.Execute("65+25=90", $true, $false, $false, $false, $false, $true, 1, $false, "99-22=77", 2)
$d.Content.Find.Execute("87+6=93", $true, $false, $false, $false, $false, $true, 1, $false, "9+39=48", 2)
$d.Content.Find.Execute("71-21=50", $true, $false, $false, $false, $false, $true, 1, $false, "50-14=36", 2)
$d.Content.Find.Execute("73+15=88", $true, $false, $false, $false, $false, $true, 1, $false, "4+61=65", 2)
$d.Content.Find.Execute("60-30=30", $true, $false, $false, $false, $false, $true, 1, $false, "75-59=16", 2)
$d.Content.Find.Execute("40+14=54", $true, $false, $false, $false, $false, $true, 1, $false, "66-46=20", 2)
$d.Content.Find.Execute("8+5=13", $true, $false, $false, $false, $false, $true, 1, $false, "46+35=81", 2)
$d.Content.Find.Execute("44+3=47", $true, $false, $false, $false, $false, $true, 1, $false, "75-42=33", 2)
$d.Content.Find.Execute("9+70=79", $true, $false, $false, $false, $false, $true, 1, $false, "22+3=25", 2)
$d.Content.Find.Execute("63-29=34", $true, $false, $false, $false, $false, $true, 1, $false, "46+44=90", 2)
$d.Content.Find.Execute("22+6=28", $true, $false, $false, $false, $false, $true, 1, $false, "96-39=57", 2)
$d.Content.Find.Execute("93+1=94", $true, $false, $false, $false, $false, $true, 1, $false, "75-19=56", 2)
$d.Content.Find.Execute("39+38=77", $true, $false, $false, $false, $false, $true, 1, $false, "27-0=27", 2)
$d.Content.Find.Execute("46-42=4", $true, $false, $false, $false, $false, $true, 1, $false, "69-66=3", 2)
$d.Content.Find.Execute("84-2=82", $true, $false, $false, $false, $false, $true, 1, $false, "55-50=5", 2)
$d.Content.Find.Execute("28-10=18", $true, $false, $false, $false, $false, $true, 1, $false, "76-71=5", 2)
$d.Content.Find.Execute("36+39=75", $true, $false, $false, $false, $false, $true, 1, $false, "25-24=1", 2)
$d.Content.Find.Execute("82-24=58", $true, $false, $false, $false, $false, $true, 1, $false, "15+22=37", 2)
$d.Content.Find.Execute("79-43=36", $true, $false, $false, $false, $false, $true, 1, $false, "46+15=61", 2)
$d.Content.Find.Execute("75-51=24", $true, $false, $false, $false, $false, $true, 1, $false, "93-64=29", 2)
$d.Content.Find.Execute("91+4=95", $true, $false, $false, $false, $false, $true, 1, $false, "65-43=22", 2)
$d.Content.Find.Execute("90-48=42", $true, $false, $false, $false, $false, $true, 1, $false, "41+6=47", 2)
$d.Content.Find.Execute("59-38=21", $true, $false, $false, $false, $false, $true, 1, $false, "78-28=50", 2)
$d.Content.Find.Execute("34+52=86", $true, $false, $false, $false, $false, $true, 1, $false, "29-25=4", 2)
$d.Content.Find.Execute("9+25=34", $true, $false, $false, $false, $false, $true, 1, $false, "77-71=6", 2)
$d.Content.Find.Execute("84-42=42", $true, $false, $false, $false, $false, $true, 1, $false, "98-52=46", 2)
$d.Content.Find.Execute("25+48=73", $true, $false, $false, $false, $false, $true, 1, $false, "74-48=26", 2)
$d.Content.Find.Execute("16+13=29", $true, $false, $false, $false, $false, $true, 1, $false, "33-2=31", 2)
$d.Content.Find.Execute("27+1=28", $true, $false, $false, $false, $false, $true, 1, $false, "70-53=17", 2)
$d.Content.Find.Execute("61-19=42", $true, $false, $false, $false, $false, $true, 1, $false, "86-22=64", 2)
$d.Content.Find.Execute("76-62=14", $true, $false, $false, $false, $false, $true, 1, $false, "90-6=84", 2)
$d.Content.Find.Execute("98+0=98", $true, $false, $false, $false, $false, $true, 1, $false, "83-23=60", 2)
$d.Content.Find.Execute("49+49=98", $true, $false, $false, $false, $false, $true, 1, $false, "89-87=2", 2)
$d.Content.Find.Execute("55-23=32", $true, $false, $false, $false, $false, $true, 1, $false, "7+89=96", 2)
$d.Content.Find.Execute("16+60=76", $true, $false, $false, $false, $false, $true, 1, $false, "87-55=32", 2)
$d.Content.Find.Execute("26+34=60", $true, $false, $false, $false, $false, $true, 1, $false, "49-4=45", 2)
$d.Content.Find.Execute("95-85=10", $true, $false, $false, $false, $false, $true, 1, $false, "80+14=94", 2)
$d.Content.Find.Execute("45+28=73", $true, $false, $false, $false, $false, $true, 1, $false, "93-79=14", 2)
$d.Content.Find.Execute("70-28=42", $true, $false, $false, $false, $false, $true, 1, $false, "12+43=55", 2)
$d.Content.Find.Execute("48-26=22", $true, $false, $false, $false, $false, $true, 1, $false, "35+28=63", 2)
$d.Content.Find.Execute("84-40=44", $true, $false, $false, $false, $false, $true, 1, $false, "38+1=39", 2)
$d.Content.Find.Execute("21+59=80", $true, $false, $false, $false, $false, $true, 1, $false, "85-74=11", 2)
$d.Content.Find.Execute("12+71=83", $true, $false, $false, $false, $false, $true, 1, $false, "99-65=34", 2)
$d.Content.Find.Execute("17+53=70", $true, $false, $false, $false, $false, $true, 1, $false, "35+35=70", 2)
$d.Content.Find.Execute("99-13=86", $true, $false, $false, $false, $false, $true, 1, $false, "71-17=54", 2)
$d.Content.Find.Execute("36+45=81", $true, $false, $false, $false, $false, $true, 1, $false, "18-16=2", 2)
$d.Content.Find.Execute("80+7=87", $true, $false, $false, $false, $false, $true, 1, $false, "48+29=77", 2)
$d.Content.Find.Execute("61-51=10", $true, $false, $false, $false, $false, $true, 1, $false, "53+24=77", 2)
$d.Content.Find.Execute("61+19=80", $true, $false, $false, $false, $false, $true, 1, $false, "77-11=66", 2)
$d.Content.Find.Execute("18-14=4", $true, $false, $false, $false, $false, $true, 1, $false, "4+87=91", 2)
$d.Content.Find.Execute("15+64=79", $true, $false, $false, $false, $false, $true, 1, $false, "31+25=56", 2)
$d.Content.Find.Execute("24+20=44", $true, $false, $false, $false, $false, $true, 1, $false, "43+33=76", 2)
